$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values look like plain numbers (e.g. "1.000",
# "27.722.45"). Force those specific cells to Text format first so
# Excel stores them verbatim instead of coercing/trimming them as
# numeric values.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D29",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '27.722.45'
$ws.Range("E2").Value = '  -1.37%  '

# Row 3
$ws.Range("D3").Value = '1.795.81'
$ws.Range("E3").Value = '  +0.15%  '

# Row 4
$ws.Range("D4").Value = '0.9981'
$ws.Range("E4").Value = '  -0.49%  '

# Row 5
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  -0.24%  '

# Row 6
$ws.Range("D6").Value = '305.91'
$ws.Range("E6").Value = '  -2.60%  '

# Row 7
$ws.Range("D7").Value = '0.4963'
$ws.Range("E7").Value = '  -4.97%  '

# Row 8
$ws.Range("D8").Value = '0.3851'
$ws.Range("E8").Value = '  +1.01%  '

# Row 9
$ws.Range("D9").Value = '0.09419'
$ws.Range("E9").Value = '  +18.47%  '

# Row 10
$ws.Range("D10").Value = '1.090'
$ws.Range("E10").Value = '  -0.18%  '

# Row 11
$ws.Range("D11").Value = '40.42'
$ws.Range("E11").Value = '  -2.41%  '

# Row 12
$ws.Range("D12").Value = '6.247'
$ws.Range("E12").Value = '  -0.41%  '

# Row 13
$ws.Range("D13").Value = '0.9982'
$ws.Range("E13").Value = '  -0.47%  '

# Row 14
$ws.Range("D14").Value = '20.40'
$ws.Range("E14").Value = '  -0.35%  '

# Row 15
$ws.Range("D15").Value = '1.795.13'
$ws.Range("E15").Value = '  -0.20%  '

# Row 16
$ws.Range("D16").Value = '7.128'
$ws.Range("E16").Value = '  -2.18%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001106'
$ws.Range("E17").Value = '  +2.07%  '

# Row 18
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '91.92'
$ws.Range("E18").Value = '  +0.11%  '

# Row 19
$ws.Range("D19").Value = '0.06536'
$ws.Range("E19").Value = '  -0.52%  '

# Row 20
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.20%  '

# Row 21
$ws.Range("D21").Value = '17.00'
$ws.Range("E21").Value = '  -1.72%  '

# Row 22
$ws.Range("D22").Value = '5.889'
$ws.Range("E22").Value = '  -0.84%  '

# Row 23
$ws.Range("D23").Value = '27.727.15'
$ws.Range("E23").Value = '  -1.48%  '

# Row 24
$ws.Range("D24").Value = '10.93'
$ws.Range("E24").Value = '  -1.74%  '

# Row 25
$ws.Range("D25").Value = '2.221'
$ws.Range("E25").Value = '  -1.92%  '

# Row 26
$ws.Range("D26").Value = '156.79'
$ws.Range("E26").Value = '  -2.37%  '

# Row 27
$ws.Range("D27").Value = '20.37'
$ws.Range("E27").Value = '  -0.20%  '

# Row 28
$ws.Range("D28").Value = '1.999.52'
$ws.Range("E28").Value = '  +0.07%  '

# Row 29
$ws.Range("D29").Value = '2.391'
$ws.Range("E29").Value = '  +2.72%  '

# Row 30
$ws.Range("E30").Value = '  +2.71%  '

# Row 31
$ws.Range("D31").Value = '0.1067'
$ws.Range("E31").Value = '  -0.98%  '

# Row 32
$ws.Range("D32").Value = '1.049'
$ws.Range("E32").Value = '  -0.11%  '

# Row 33
$ws.Range("D33").Value = '3.601'
$ws.Range("E33").Value = '  -1.96%  '

# Row 34
$ws.Range("D34").Value = '5.501'
$ws.Range("E34").Value = '  -0.45%  '

# Row 35
$ws.Range("D35").Value = '0.06789'
$ws.Range("E35").Value = '  -5.95%  '

# Row 36
$ws.Range("D36").Value = '8.836'
$ws.Range("E36").Value = '  +0.53%  '

# Row 37
$ws.Range("D37").Value = '0.02294'
$ws.Range("E37").Value = '  -0.98%  '

# Row 38
$ws.Range("D38").Value = '0.2122'
$ws.Range("E38").Value = '  -0.80%  '

# Row 39
$ws.Range("D39").Value = '11.34'
$ws.Range("E39").Value = '  -7.02%  '

# Row 40
$ws.Range("D40").Value = '4.895'
$ws.Range("E40").Value = '  -3.24%  '

# Row 41
$ws.Range("D41").Value = '0.6107'
$ws.Range("E41").Value = '  -0.49%  '

# Row 42
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.14%  '

# Row 43
$ws.Range("D43").Value = '1.136'
$ws.Range("E43").Value = '  -2.32%  '

# Row 44
$ws.Range("D44").Value = '12.89'
$ws.Range("E44").Value = '  -2.59%  '

# Row 45
$ws.Range("D45").Value = '0.5845'
$ws.Range("E45").Value = '  -2.04%  '

# Row 46
$ws.Range("D46").Value = '3.664'
$ws.Range("E46").Value = '  -2.86%  '

# Row 47
$ws.Range("D47").Value = '1.264'
$ws.Range("E47").Value = '  -7.38%  '

# Row 48
$ws.Range("D48").Value = '123.02'
$ws.Range("E48").Value = '  -3.59%  '

# Row 49
$ws.Range("D49").Value = '1.924'
$ws.Range("E49").Value = '  +0.57%  '

# Row 50
$ws.Range("D50").Value = '1.166'
$ws.Range("E50").Value = '  -5.17%  '

# Row 51
$ws.Range("D51").Value = '0.06678'
$ws.Range("E51").Value = '  -0.88%  '
